$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Request")

# The Post-Consolidation section had a "Person State Identification ID" row
# (based on nc:PersonStateIdentification) that duplicated the Pre-Consolidation
# row. Per the commit, "PersonStateIdentification" is removed from the
# consolidation/expungement/identifier-update mapping altogether, so delete
# the Post-Consolidation "Person State Identification ID" row (row 8). This
# shifts the "Person State Fingerprint ID" / "Person FBI Identification ID"
# rows up from 9/10 to 8/9.
$ws.Rows("8:8").Delete() | Out-Null

# The Pre-Consolidation section's "Person State Identification ID" row (row 5)
# is likewise replaced by a "Person State Fingerprint ID" row, matching the
# pattern now used consistently in the Post-Consolidation section.
$ws.Range("A5").Value2 = "Person State Fingerprint ID"
$ws.Range("B5").Value2 = "An identification of a person based on a Fingerprint ID."
$ws.Range("C5").Value2 = "chc-report-doc:CriminalHistoryConsolidationReport/nc:Person/chc-report-ext:PreConsolidationIdentifiers/j:PersonStateFingerprintIdentification/nc:IdentificationID"

# Update the remaining NEIM 3.0 Mapping (column C) entries so the XPath-like
# mapping strings use the "chc-report-doc"/"chc-report-ext" namespace prefixes
# (replacing the old "/CHcr-doc"/"CHcr-ext" prefixes, and dropping the
# leading "/").
$ws.Range("C6").Value2 = "chc-report-doc:CriminalHistoryConsolidationReport/nc:Person/chc-report-ext:PreConsolidationIdentifiers/j:PersonFBIIdentification/nc:IdentificationID"
$ws.Range("C8").Value2 = "chc-report-doc:CriminalHistoryConsolidationReport/nc:Person/chc-report-ext:PostConsolidationIdentifiers/j:PersonStateFingerprintIdentification/nc:IdentificationID"
$ws.Range("C9").Value2 = "chc-report-doc:CriminalHistoryConsolidationReport/nc:Person/chc-report-ext:PostConsolidationIdentifiers/j:PersonFBIIdentification/nc:IdentificationID"

# Update the selection to reflect where the user last clicked after the edit.
$ws.Range("C9").Select() | Out-Null
